# 4th Stab - Added Marketeat
# Rename the "data" sheet to "2018" and move the selection to F43.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab (was "data" -> "2018")
$ws.Name = "2018"

# Make sure it's the active sheet, then move the active selection
$ws.Activate()
$ws.Range("F43").Select()
